$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'334.20"
$ws.Range("E2").Value = "'1.65%"
$ws.Range("D3").Value = "'43.84"
$ws.Range("E3").Value = "'6.02%"
$ws.Range("D4").Value = "'5.815"
$ws.Range("E4").Value = "'3.07%"
$ws.Range("E5").Value = "'1.99%"
$ws.Range("D6").Value = "'8.802"
$ws.Range("E6").Value = "'0.62%"
$ws.Range("D7").Value = "'1.991"
$ws.Range("E7").Value = "'-1.08%"
$ws.Range("D9").Value = "'0.9412"
$ws.Range("E9").Value = "'2.12%"
$ws.Range("E10").Value = "'-2.20%"
$ws.Range("D11").Value = "'0.1949"
$ws.Range("E11").Value = "'-0.65%"
$ws.Range("D12").Value = "'0.09679"
$ws.Range("E12").Value = "'2.52%"
$ws.Range("D13").Value = "'0.04608"
$ws.Range("E13").Value = "'23.05%"
$ws.Range("D14").Value = "'0.1068"
$ws.Range("E14").Value = "'0.82%"
$ws.Range("D15").Value = "'0.001298"
$ws.Range("E15").Value = "'-0.19%"
$ws.Range("D16").Value = "'0.005947"
$ws.Range("E16").Value = "'-2.83%"
$ws.Range("D17").Value = "'3.498"
$ws.Range("E17").Value = "'1.52%"
$ws.Range("D18").Value = "'4.506"
$ws.Range("E18").Value = "'0.22%"
$ws.Range("D20").Value = "'8.753"
$ws.Range("E20").Value = "'4.67%"
$ws.Range("E21").Value = "'-0.70%"
$ws.Range("D23").Value = "'0.04406"
$ws.Range("E23").Value = "'0.28%"
$ws.Range("E24").Value = "'-0.01%"
$ws.Range("D25").Value = "'0.004411"
$ws.Range("E25").Value = "'2.13%"
$ws.Range("D26").Value = "'0.0001271"
$ws.Range("E26").Value = "'5.75%"
$ws.Range("D27").Value = "'0.0003995"
$ws.Range("E39").Value = "'0.24%"
$ws.Range("D40").Value = "'0.05721"
$ws.Range("E40").Value = "'5.93%"
$ws.Range("D41").Value = "'0.007961"
$ws.Range("E41").Value = "'3.63%"
$ws.Range("E42").Value = "'0.85%"
$ws.Range("D43").Value = "'0.009044"
$ws.Range("E43").Value = "'1.03%"
$ws.Range("E44").Value = "'-3.28%"
$ws.Range("D45").Value = "'0.01052"
$ws.Range("E45").Value = "'-8.81%"
$ws.Range("D46").Value = "'0.00007217"
$ws.Range("E46").Value = "'9.28%"
$ws.Range("E47").Value = "'-0.06%"
$ws.Range("D48").Value = "'0.003240"
$ws.Range("E48").Value = "'1.49%"
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("E51").Value = "'-0.06%"
